$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Snapshot the current (pre-edit) time values in A3:A13 and copy them,
#     in order, into a new helper column D (matching formatting) ---
for ($i = 3; $i -le 13; $i++) {
    $ws.Cells.Item($i, 4).Value = $ws.Cells.Item($i, 1).Value2
    $ws.Cells.Item($i, 4).NumberFormat = "h:mm"
}

# --- Sort the helper column D ascending, recording a sortState ---
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("D3"))
$ws.Sort.SetRange($ws.Range("D3:D13"))
$ws.Sort.Header = 2
$ws.Sort.Apply()

# --- Fix up a data-entry typo in A12 (0:30:46 -> 0:46:40 style correction) ---
$ws.Range("A12").Value = 0.644444444444444

# --- Label the summary statistics in column B ---
$ws.Range("B15").Value = "mean"
$ws.Range("B16").Value = "median"
$ws.Range("B17").Value = "mode"

# --- Drop the standard-deviation summary row; keep the (now-empty) cell ---
$ws.Range("A18").ClearContents()

# --- Restore the selection to the cell that was just edited ---
$ws.Range("A12").Select()
